$d = $word.ActiveDocument

$replacements = @(
    @{old="856÷5=171, 1"; new="209÷9=23, 2"},
    @{old="156÷5=31, 1"; new="154÷4=38, 2"},
    @{old="154÷3=51, 1"; new="616÷6=102, 4"},
    @{old="775÷8=96, 7"; new="244÷3=81, 1"},
    @{old="954÷3=318, 0"; new="820÷8=102, 4"},
    @{old="155÷6=25, 5"; new="194÷8=24, 2"},
    @{old="484÷3=161, 1"; new="272÷5=54, 2"},
    @{old="901÷4=225, 1"; new="654÷9=72, 6"},
    @{old="817÷5=163, 2"; new="374÷7=53, 3"},
    @{old="110÷7=15, 5"; new="155÷7=22, 1"},
    @{old="229÷5=45, 4"; new="487÷8=60, 7"},
    @{old="920÷5=184, 0"; new="370÷6=61, 4"},
    @{old="649÷7=92, 5"; new="939÷2=469, 1"},
    @{old="384÷8=48, 0"; new="798÷2=399, 0"},
    @{old="235÷8=29, 3"; new="666÷7=95, 1"},
    @{old="971÷9=107, 8"; new="284÷7=40, 4"},
    @{old="557÷4=139, 1"; new="656÷5=131, 1"},
    @{old="951÷7=135, 6"; new="104÷3=34, 2"},
    @{old="408÷6=68, 0"; new="799÷9=88, 7"},
    @{old="688÷9=76, 4"; new="667÷7=95, 2"},
    @{old="428÷9=47, 5"; new="574÷6=95, 4"},
    @{old="943÷4=235, 3"; new="742÷4=185, 2"},
    @{old="861÷8=107, 5"; new="611÷9=67, 8"},
    @{old="543÷9=60, 3"; new="975÷8=121, 7"},
    @{old="274÷4=68, 2"; new="224÷4=56, 0"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
